# Updates the cryptos price list (columns D and E) for rows 2-51
# to the latest scraped values, matching the GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.491.42'
$ws.Range("E2").Value = '  -0.45%  '
$ws.Range("D3").Value = '1.818.86'
$ws.Range("E3").Value = '  -0.72%  '
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = "'316.34"
$ws.Range("E5").Value = '  -0.50%  '
$ws.Range("E6").Value = '  +0.14%  '
$ws.Range("D7").Value = "'0.5146"
$ws.Range("E7").Value = '  -3.74%  '
$ws.Range("D8").Value = "'0.3863"
$ws.Range("E8").Value = '  -3.28%  '
$ws.Range("D9").Value = "'0.08429"
$ws.Range("E9").Value = '  +9.04%  '
$ws.Range("E10").Value = '  -0.15%  '
$ws.Range("D11").Value = "'1.109"
$ws.Range("E11").Value = '  -1.02%  '
$ws.Range("E12").Value = '  +1.36%  '
$ws.Range("D13").Value = "'21.03"
$ws.Range("E13").Value = '  +0.16%  '
$ws.Range("E14").Value = '  +0.17%  '
$ws.Range("E15").Value = '  -1.47%  '
$ws.Range("D16").Value = '1.817.22'
$ws.Range("E16").Value = '  -0.75%  '
$ws.Range("E17").Value = '  +4.43%  '
$ws.Range("D18").Value = "'93.15"
$ws.Range("E18").Value = '  -0.43%  '
$ws.Range("D19").Value = "'0.06713"
$ws.Range("E19").Value = '  +2.17%  '
$ws.Range("E20").Value = '  -0.37%  '
$ws.Range("E21").Value = '  +0.10%  '
$ws.Range("E22").Value = '  -0.17%  '
$ws.Range("D23").Value = '28.518.75'
$ws.Range("E23").Value = '  -0.40%  '
$ws.Range("E24").Value = '  +1.73%  '
$ws.Range("D25").Value = "'2.272"
$ws.Range("E25").Value = '  +1.38%  '
$ws.Range("D26").Value = "'21.21"
$ws.Range("E26").Value = '  +2.17%  '
$ws.Range("D27").Value = "'158.64"
$ws.Range("E27").Value = '  +1.49%  '
$ws.Range("D28").Value = '2.027.13'
$ws.Range("E28").Value = '  -0.74%  '
$ws.Range("D29").Value = "'2.405"
$ws.Range("E29").Value = '  -0.66%  '
$ws.Range("D30").Value = "'125.86"
$ws.Range("E30").Value = '  +0.55%  '
$ws.Range("D31").Value = "'1.094"
$ws.Range("E31").Value = '  -4.11%  '
$ws.Range("D32").Value = "'0.1080"
$ws.Range("E32").Value = '  -3.54%  '
$ws.Range("D33").Value = "'5.746"
$ws.Range("E33").Value = '  -0.06%  '
$ws.Range("D34").Value = "'3.690"
$ws.Range("E34").Value = '  +1.04%  '
$ws.Range("D35").Value = "'0.07375"
$ws.Range("E35").Value = '  +1.71%  '
$ws.Range("E36").Value = '  -1.34%  '
$ws.Range("D37").Value = "'0.02363"
$ws.Range("E37").Value = '  +0.54%  '
$ws.Range("D38").Value = "'5.224"
$ws.Range("D39").Value = "'8.822"
$ws.Range("E39").Value = '  -0.87%  '
$ws.Range("D40").Value = "'0.6319"
$ws.Range("E40").Value = '  +0.09%  '
$ws.Range("D41").Value = "'11.24"
$ws.Range("E41").Value = '  -1.45%  '
$ws.Range("D42").Value = "'1.194"
$ws.Range("E42").Value = '  -0.33%  '
$ws.Range("D43").Value = "'1.402"
$ws.Range("E43").Value = '  +0.48%  '
$ws.Range("D44").Value = "'13.48"
$ws.Range("E44").Value = '  -0.80%  '
$ws.Range("D45").Value = "'3.760"
$ws.Range("E45").Value = '  +1.08%  '
$ws.Range("D46").Value = "'0.5899"
$ws.Range("E46").Value = '  -0.21%  '
$ws.Range("D47").Value = "'125.79"
$ws.Range("E47").Value = '  +0.46%  '
$ws.Range("E48").Value = '  -0.65%  '
$ws.Range("D49").Value = "'1.195"
$ws.Range("E49").Value = '  -0.17%  '
$ws.Range("D50").Value = "'0.06981"
$ws.Range("E50").Value = '  +0.56%  '
$ws.Range("D51").Value = "'74.08"
$ws.Range("E51").Value = '  -0.65%  '
